# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1339
$ws1.Range("F4").Value  = 14682
$ws1.Range("F5").Value  = 17934
$ws1.Range("F6").Value  = 17934
$ws1.Range("F8").Value  = 80
$ws1.Range("F10").Value = 221
$ws1.Range("F16").Value = 49
$ws1.Range("F17").Value = 58
$ws1.Range("F18").Value = 164
$ws1.Range("F20").Value = 1350
$ws1.Range("F21").Value = 150
$ws1.Range("F22").Value = 78
$ws1.Range("F24").Value = 217
$ws1.Range("F25").Value = 7406
$ws1.Range("F26").Value = 982
$ws1.Range("F28").Value = 46
$ws1.Range("F31").Value = 5886
$ws1.Range("F34").Value = 148
$ws1.Range("F36").Value = 235
$ws1.Range("F37").Value = 5149
$ws1.Range("F39").Value = 35

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1339
$ws4.Range("F4").Value  = 14682
$ws4.Range("F5").Value  = 17934
$ws4.Range("F6").Value  = 17934
$ws4.Range("F8").Value  = 80
$ws4.Range("F10").Value = 221
$ws4.Range("F16").Value = 49
$ws4.Range("F17").Value = 58
$ws4.Range("F18").Value = 164
$ws4.Range("F20").Value = 1350
$ws4.Range("F21").Value = 150
$ws4.Range("F22").Value = 78
$ws4.Range("F25").Value = 217
$ws4.Range("F26").Value = 7406
$ws4.Range("F27").Value = 982
$ws4.Range("F29").Value = 46
$ws4.Range("F33").Value = 5886
$ws4.Range("F36").Value = 148
$ws4.Range("F38").Value = 235
$ws4.Range("F39").Value = 5149
$ws4.Range("F41").Value = 35

$wb.Save()
